# HonsProgress.xlsx - PDD sheet progress update
# "pdd section 3 + prep for meeting"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PDD")

# --- Constraints (row 7) no longer a blocking issue -> downgrade from the
#     "Bad" red/pink flag to the lighter red-text/boxed flag (same format
#     already used on row 13, "Assumptions"). Copy that cell's formatting so
#     we land on the exact same style instead of inventing a new one.
$ws.Range("D13").Copy()
$ws.Range("D7").PasteSpecial(-4122)

# --- Project Eval (row 18) becomes the new outstanding item -> flag it with
#     the "Bad" style, matching how rows 7/11 were flagged before this edit.
#     Grab that formatting from D11 first, before D11 itself is cleared below.
$ws.Range("D11").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Section 2 draft completion bumped up from 70% to 75% ---
$ws.Range("F10").Value = 0.75   # Project Operation
$ws.Range("F11").Value = 0.75   # Options
$ws.Range("F12").Value = 0.75   # Risks Mits

# Options (row 11) is no longer flagged -> clear back to the plain style
# used elsewhere in the table (copy format from row 12, "Risks Mits").
$ws.Range("D12").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Section 3 / Project Control (row 17) progressed significantly ---
$ws.Range("E17").Value = 144
$ws.Range("F17").Value = 0.75

# Leave the cursor where it ended up after making these edits.
$ws.Range("F26").Select() | Out-Null
